$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176317930221558
$ws.Range("B1").Value = 2.145128726959229
$ws.Range("C1").Value = 3.036573886871338
$ws.Range("D1").Value = 3.549987077713013
$ws.Range("E1").Value = 1.570531606674194
